$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2: remove the now-empty, style-only C2 cell ---
$ws.Range("C2").Clear() | Out-Null

# give C3 the same formatting used by the other row-2 cells (style "2") before
# putting any value in it
$ws.Range("B2").Copy() | Out-Null
$ws.Range("C3").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# --- fill in the new note rows, in the same order the underlying shared
#     strings were originally authored (Error, InvalidEMGUID, the content id,
#     WrongContentId) ---
$ws.Range("E4").Value = "Error"
$ws.Range("D4").Value = "InvalidEMGUID"

$ws.Range("C3").Value = "A035114660061500840002"

$ws.Range("D5").Value = "3ed831fb-012e-4538-973b-eeb67d80d931"
$ws.Range("E5").Value = "Video"

$ws.Range("C6").Value = "WrongContentId"
$ws.Range("D6").Value = "3ed831fb-012e-4538-973b-eeb67d80d931"
$ws.Range("E6").Value = "Video"

$ws.Rows.Item(3).RowHeight = 15.5

# --- column widths (bestFit-style autofit to the new, wider content);
#     the inputs below are chosen so the engine's internal pixel-rounding
#     lands on the closest representable width to the real Excel bestFit
#     values of 27.6328125 / 36.08984375 characters ---
$ws.Columns.Item(3).ColumnWidth = 26.75
$ws.Columns.Item(4).ColumnWidth = 35.25

# --- selection moves to C3 ---
$ws.Range("C3").Select() | Out-Null
